$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 3599.5293
$ws.Range("I100").Value = 2286.9092
$ws.Range("J100").Value = 6006
$ws.Range("K100").Value = 2286.9092
$ws.Range("L100").Value = 6006
$ws.Range("M100").Value = -1745.9092
$ws.Range("N100").Value = -7088
$ws.Range("I106").Value = 3937.5
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 3937.5
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -3306.5
$ws.Range("N106").ClearContents()
$ws.Range("H110").Value = 75000
$ws.Range("J110").Value = 75000
$ws.Range("L110").Value = 75000
$ws.Range("N110").Value = -83180
$ws.Range("H129").Value = 2932.1316
$ws.Range("J129").Value = 2992.3333
$ws.Range("L129").Value = 8976.999899999999
$ws.Range("N129").Value = -18976.9999
$ws.Range("H131").Value = 6298.5835
$ws.Range("I131").Value = 722.5
$ws.Range("J131").Value = 9086.625
$ws.Range("K131").Value = 2167.5
$ws.Range("L131").Value = 27259.875
$ws.Range("M131").Value = 2872.5
$ws.Range("N131").Value = -37339.875
$ws.Range("H138").Value = 4063.5098
$ws.Range("I138").Value = 1368.1
$ws.Range("J138").Value = 4720.927
$ws.Range("K138").Value = 4104.299999999999
$ws.Range("L138").Value = 14162.781
$ws.Range("M138").Value = 1035.700000000001
$ws.Range("N138").Value = -24442.781
$ws.Range("H141").Value = 2033.4615
$ws.Range("I141").Value = 1891.3235
$ws.Range("J141").Value = 3000
$ws.Range("K141").Value = 5673.970499999999
$ws.Range("L141").Value = 9000
$ws.Range("M141").Value = -493.9704999999994
$ws.Range("N141").Value = -19360
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3435.7344
$ws.Range("I32").Value = 2349.611
$ws.Range("K32").Value = 2349.611
$ws.Range("M32").Value = -2062.611
$ws.Range("H45").Value = 2356.1667
$ws.Range("I45").Value = 1141.5555
$ws.Range("K45").Value = 1141.5555
$ws.Range("M45").Value = -764.5554999999999
$ws.Range("H97").Value = 944
$ws.Range("I97").Value = 754.75
$ws.Range("K97").Value = 754.75
$ws.Range("M97").Value = -258.75
$ws.Range("H132").Value = 2972.9673
$ws.Range("I132").Value = 3372.1836
$ws.Range("J132").Value = 1342.8334
$ws.Range("K132").Value = 10116.5508
$ws.Range("L132").Value = 4028.5002
$ws.Range("M132").Value = -7586.550799999999
$ws.Range("N132").Value = -9088.5002
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 2822.4614
$ws.Range("I5").Value = 614.5
$ws.Range("J5").Value = 4715
$ws.Range("K5").Value = 614.5
$ws.Range("L5").Value = 4715
$ws.Range("M5").Value = -501.5
$ws.Range("N5").Value = -4941
$ws.Range("H94").Value = 1078.3529
$ws.Range("J94").Value = 910.6667
$ws.Range("L94").Value = 910.6667
$ws.Range("N94").Value = -1812.6667
$ws.Range("H99").Value = 2675.087
$ws.Range("I99").Value = 2151.3572
$ws.Range("K99").Value = 2151.3572
$ws.Range("M99").Value = -653.3571999999999
$ws.Range("H105").Value = 614.9
$ws.Range("I105").Value = 576.44446
$ws.Range("K105").Value = 576.44446
$ws.Range("M105").Value = 1170.55554
$ws.Range("H134").Value = 2901.6365
$ws.Range("I134").Value = 2559.7144
$ws.Range("K134").Value = 7679.1432
$ws.Range("M134").Value = -5144.1432
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1781.421
$ws.Range("I132").Value = 1802.6111
$ws.Range("K132").Value = 5407.8333
$ws.Range("M132").Value = -2877.8333
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 14494918
$ws.Range("J4").Value = 46000028
$ws.Range("L4").Value = 138000084
$ws.Range("N4").Value = -138000308
$ws.Range("H7").Value = 912.6
$ws.Range("I7").Value = 912.6
$ws.Range("K7").Value = 2737.8
$ws.Range("M7").Value = -2625.8
$ws.Range("H80").Value = 2249.25
$ws.Range("J80").Value = 2249.25
$ws.Range("L80").Value = 6747.75
$ws.Range("N80").Value = -8619.75
$ws.Range("H83").Value = 2249.25
$ws.Range("J83").Value = 2249.25
$ws.Range("L83").Value = 20243.25
$ws.Range("N83").Value = -29603.25
$ws.Range("H123").Value = 1117.5555
$ws.Range("I123").Value = 1117.5555
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 3352.6665
$ws.Range("L123").Value = 0
$ws.Range("M123").Value = -902.6664999999998
$ws.Range("N123").ClearContents()
$ws.Range("H129").Value = 3268.7856
$ws.Range("I129").Value = 643.3333
$ws.Range("K129").Value = 1929.9999
$ws.Range("M129").Value = 3070.0001
$ws.Range("H131").Value = 1834.56
$ws.Range("I131").Value = 806.4
$ws.Range("J131").Value = 2520
$ws.Range("K131").Value = 2419.2
$ws.Range("L131").Value = 7560
$ws.Range("M131").Value = 2620.8
$ws.Range("N131").Value = -17640
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 18000
$ws.Range("J39").Value = 18000
$ws.Range("L39").Value = 18000
$ws.Range("N39").Value = -19064
$ws.Range("H49").Value = 24999.334
$ws.Range("J49").Value = 24999.334
$ws.Range("L49").Value = 24999.334
$ws.Range("N49").Value = -25367.334
$ws.Range("H52").Value = 36857.285
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 36857.285
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 36857.285
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -37375.285
$ws.Range("H122").Value = 3324.5
$ws.Range("I122").Value = 2759.4
$ws.Range("J122").Value = 3728.1428
$ws.Range("K122").Value = 8278.200000000001
$ws.Range("L122").Value = 11184.4284
$ws.Range("M122").Value = -5828.200000000001
$ws.Range("N122").Value = -16084.4284
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 8500
$ws.Range("I5").Value = 7000
$ws.Range("J5").Value = 10000
$ws.Range("K5").Value = 7000
$ws.Range("L5").Value = 10000
$ws.Range("M5").Value = -6887
$ws.Range("N5").Value = -10226
$ws.Range("H40").Value = 5659.1333
$ws.Range("I40").Value = 3200
$ws.Range("K40").Value = 3200
$ws.Range("M40").Value = -3064
$ws.Range("H42").Value = 1000000
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H49").Value = 1000000
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H100").Value = 7000
$ws.Range("I100").Value = 5500
$ws.Range("J100").Value = 7600
$ws.Range("K100").Value = 5500
$ws.Range("L100").Value = 7600
$ws.Range("M100").Value = -4959
$ws.Range("N100").Value = -8682
$ws.Range("H122").Value = 2466.6667
$ws.Range("I122").Value = 2456.1538
$ws.Range("J122").Value = 2535
$ws.Range("K122").Value = 7368.4614
$ws.Range("L122").Value = 7605
$ws.Range("M122").Value = -4918.4614
$ws.Range("N122").Value = -12505
$ws.Range("H132").Value = 2393
$ws.Range("I132").Value = 1743.129
$ws.Range("J132").Value = 3400.3
$ws.Range("K132").Value = 5229.387
$ws.Range("L132").Value = 10200.9
$ws.Range("M132").Value = -2699.387
$ws.Range("N132").Value = -15260.9
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 1767.5
$ws.Range("I14").Value = 2291.5
$ws.Range("J14").Value = 1374.5
$ws.Range("K14").Value = 2291.5
$ws.Range("L14").Value = 1374.5
$ws.Range("M14").Value = -2123.5
$ws.Range("N14").Value = -1710.5
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("N21").ClearContents()
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("N35").ClearContents()
$ws.Range("H38").Value = 24249.25
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 24249.25
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 24249.25
$ws.Range("M38").ClearContents()
$ws.Range("N38").Value = -25195.25
$ws.Range("H96").Value = 9400.6
$ws.Range("I96").Value = 3499.5
$ws.Range("J96").Value = 13334.667
$ws.Range("K96").Value = 3499.5
$ws.Range("L96").Value = 13334.667
$ws.Range("M96").Value = -2126.5
$ws.Range("N96").Value = -16080.667
$ws.Range("H132").Value = 6758.6743
$ws.Range("I132").Value = 6871.684
$ws.Range("K132").Value = 20615.052
$ws.Range("M132").Value = -18085.052
$ws.Range("H136").Value = 1848.509
$ws.Range("I136").Value = 1193.4286
$ws.Range("J136").Value = 2994.9
$ws.Range("K136").Value = 3580.2858
$ws.Range("L136").Value = 8984.700000000001
$ws.Range("M136").Value = -1030.2858
$ws.Range("N136").Value = -14084.7